# Update JLP_2024-12.xlsx: fix transfer/receita rows and remove extra rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 12 (no longer needed)
$ws.Range("A7:F12").EntireRow.Delete()

# Row 2: Despesa / TARIFAS
$ws.Range("A2").Value = "Despesa"
$ws.Range("B2").Value = "TARIFAS"
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = "20/02/2025"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = " "

# Row 3: Despesa / MOBILIÁRIO
$ws.Range("A3").Value = "Despesa"
$ws.Range("B3").Value = "MOBILIÁRIO"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = "20/02/2025"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = " "

# Row 4: Despesa / ESCRITÓRIO
$ws.Range("A4").Value = "Despesa"
$ws.Range("B4").Value = "ESCRITÓRIO"
$ws.Range("C4").Value = 250
$ws.Range("D4").Value = "20/02/2025"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = " "

# Row 5: Receita / APLICAÇÕES FINANCEIRAS (Inquilino cleared)
$ws.Range("B5").Value = "APLICAÇÕES FINANCEIRAS"
$ws.Range("C5").Value = 250
$ws.Range("D5").Value = "20/02/2025"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = ""

# Row 6: Receita / APLICAÇÕES FINANCEIRAS (Inquilino cleared)
$ws.Range("B6").Value = "APLICAÇÕES FINANCEIRAS"
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = "20/02/2025"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = ""
